# Atualização de bases das ligas, do dia: 12-04-2024 às 20:28
# Swap the match-data (all columns except the id column A) between
# row 11 <-> row 12, and between row 83 <-> row 84.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @(11, 12),
    @(83, 84)
)

foreach ($pair in $pairs) {
    $row1 = $pair[0]
    $row2 = $pair[1]

    for ($c = 2; $c -le 29; $c++) {
        $cell1 = $ws.Cells.Item($row1, $c)
        $cell2 = $ws.Cells.Item($row2, $c)

        $v1 = $cell1.Value2
        $v2 = $cell2.Value2

        $cell1.Value = $v2
        $cell2.Value = $v1
    }
}
